# Apply updated cryptocurrency prices and volume percentages
# (values provided with a leading quote where needed so Excel stores
#  them as text rather than re-parsing as numbers, matching the
#  original inline-string cell content)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.270.14"
$ws.Range("D3").Value = "1.862.70"
$ws.Range("E3").Value = "  -1.98%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'319.24"
$ws.Range("E5").Value = "  -1.52%  "
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").Value = "'0.4379"
$ws.Range("E7").Value = "  -4.66%  "
$ws.Range("D8").Value = "'0.3718"
$ws.Range("E8").Value = "  -2.40%  "
$ws.Range("D9").Value = "'0.07529"
$ws.Range("E9").Value = "  -2.34%  "
$ws.Range("D10").Value = "'0.9388"
$ws.Range("E10").Value = "  -3.55%  "
$ws.Range("D11").Value = "'21.31"
$ws.Range("E11").Value = "  -2.85%  "
$ws.Range("D12").Value = "1.873.13"
$ws.Range("E12").Value = "  -3.17%  "
$ws.Range("D13").Value = "'6.733"
$ws.Range("E13").Value = "  -2.69%  "
$ws.Range("D14").Value = "'5.475"
$ws.Range("E14").Value = "  -3.02%  "
$ws.Range("D15").Value = "'0.06855"
$ws.Range("E15").Value = "  -3.03%  "
$ws.Range("D16").Value = "'1.005"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").Value = "'82.15"
$ws.Range("E17").Value = "  -1.97%  "
$ws.Range("D18").Value = "'0.000009106"
$ws.Range("E18").Value = "  -3.80%  "
$ws.Range("D19").Value = "'1.003"
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").Value = "'16.01"
$ws.Range("E20").Value = "  -3.60%  "
$ws.Range("D21").Value = "28.260.48"
$ws.Range("E21").Value = "  -2.24%  "
$ws.Range("D22").Value = "'5.155"
$ws.Range("E22").Value = "  -2.19%  "
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("D24").Value = "2.090.43"
$ws.Range("E24").Value = "  -2.53%  "
$ws.Range("D25").Value = "'2.038"
$ws.Range("E25").Value = "  -2.90%  "
$ws.Range("D26").Value = "'154.77"
$ws.Range("E26").Value = "  -2.14%  "
$ws.Range("D27").Value = "'18.39"
$ws.Range("E27").Value = "  -3.33%  "
$ws.Range("D28").Value = "'5.357"
$ws.Range("E28").Value = "  -4.51%  "
$ws.Range("D29").Value = "'114.34"
$ws.Range("E29").Value = "  -2.68%  "
$ws.Range("D30").Value = "'1.733"
$ws.Range("E30").Value = "  -5.87%  "
$ws.Range("D31").Value = "'0.09062"
$ws.Range("E31").Value = "  -1.99%  "
$ws.Range("D32").Value = "'0.8036"
$ws.Range("E32").Value = "  -6.02%  "
$ws.Range("D33").Value = "'4.854"
$ws.Range("E33").Value = "  -4.31%  "
$ws.Range("D34").Value = "'1.173"
$ws.Range("E34").Value = "  -5.00%  "
$ws.Range("D35").Value = "'2.944"
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("D37").Value = "'1.127"
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("D38").Value = "'0.05468"
$ws.Range("E38").Value = "  -3.51%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01951"
$ws.Range("E39").Value = "  -3.79%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.980"
$ws.Range("E40").Value = "  +8.05%  "
$ws.Range("D41").Value = "'7.150"
$ws.Range("E41").Value = "  -3.06%  "
$ws.Range("D42").Value = "'0.5256"
$ws.Range("E42").Value = "  -3.90%  "
$ws.Range("D43").Value = "'0.1675"
$ws.Range("E43").Value = "  -4.26%  "
$ws.Range("D44").Value = "'8.788"
$ws.Range("E44").Value = "  -5.08%  "
$ws.Range("D45").Value = "'2.080"
$ws.Range("E45").Value = "  +1.20%  "
$ws.Range("D46").Value = "'0.06775"
$ws.Range("E46").Value = "  -0.68%  "
$ws.Range("D47").Value = "'0.4880"
$ws.Range("D48").Value = "'0.000002533"
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("D49").Value = "'10.59"
$ws.Range("E49").Value = "  -5.27%  "
$ws.Range("D50").Value = "'107.81"
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("D51").Value = "'1.683"
$ws.Range("E51").Value = "  -4.68%  "